$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 131078252
$ws.Range("B6").Value = 57064
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 102612
$ws.Range("F6").Value = "Järpe"
$ws.Range("G6").Value = "Tetrastes bonasia"
$ws.Range("H6").Value = "(Linnaeus, 1758)"

$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("I6").Style = "Normal"

$ws.Range("K6").Value = "adult"
$ws.Range("L6").Value = "hane"
$ws.Range("M6").Value = "upprörd, varnande"
$ws.Range("N6").Value = "observerad"
$ws.Range("P6").Value = "Långmyran, Dlr"
$ws.Range("Q6").Value = 504140
$ws.Range("R6").Value = 6691587
$ws.Range("S6").Value = 20
$ws.Range("T6").Value = "Dalarna"
$ws.Range("U6").Value = "Ludvika"
$ws.Range("V6").Value = "Dalarna"
$ws.Range("W6").Value = "Grangärde"

$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2026-02-08"
$ws.Range("Y6").Style = "Normal"

$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2026-02-08"
$ws.Range("AA6").Style = "Normal"

$ws.Range("AC6").Value = "Jag lockar fram tuppen med järppipa efter att ha stött den under skidåkning."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false

$ws.Range("AT6").Formula = "'"
$ws.Range("AT6").Style = "Normal"

$ws.Range("AY6").Formula = "'"
$ws.Range("AY6").Style = "Normal"

$ws.Range("AW6").Value = "Tobias Hellberg"
$ws.Range("AX6").Value = "Tobias Hellberg"
